$wb = $excel.ActiveWorkbook

# New row-2 values coming from the server re-run of the "10 hot summers" cost
# model. Column order on every sheet is:
# A=eb B=hp C=st D=wi E=gt F=dgt G=ieh H=chp I=ac J=ab_ct K=ab_hp L=cp_ct
# M=cp_hp N=ttes O=ites
# Only columns A,B,E,G,I,L,M,N,O change; C,D,F,H,J,K stay 0 on every sheet.
# NOTE: sheet names ("2025", "2030", ...) look numeric, so they are kept out
# of hashtable keys (PowerShell would silently coerce "2025" -> [int]2025
# there) and instead carried as explicit [string] array entries.

$sheetNames = [string[]]@("2025", "2030", "2035", "2040", "2045", "2050")

$sheetValues = @(
    @{
        "A2" = 0
        "B2" = 281.6081816490545
        "E2" = 29044.8710791117
        "G2" = 8095.92571266189
        "I2" = 14901.67422553115
        "L2" = 50875.914042756
        "M2" = 11229.470312225
        "N2" = 7247.308947314462
        "O2" = 6890.471638847674
    },
    @{
        "A2" = 0
        "B2" = 3636.073638914573
        "E2" = 45890.08104099892
        "G2" = 8095.92571266189
        "I2" = 36361.75759182434
        "L2" = 63083.39540355118
        "M2" = 17442.4765456975
        "N2" = 9351.077465916849
        "O2" = 8472.471324946595
    },
    @{
        "A2" = 2148.160171297619
        "B2" = 5814.069708004545
        "E2" = 57548.65702072511
        "G2" = 8095.92571266189
        "I2" = 47220.76837855846
        "L2" = 63083.39540355118
        "M2" = 23203.20223601122
        "N2" = 13680.12751419989
        "O2" = 13142.91637740311
    },
    @{
        "A2" = 2148.160171297619
        "B2" = 5814.069708004545
        "E2" = 57548.65702072511
        "G2" = 8095.92571266189
        "I2" = 47220.76837855846
        "L2" = 63083.39540355118
        "M2" = 23203.20223601122
        "N2" = 13680.12751419989
        "O2" = 13142.91637740311
    },
    @{
        "A2" = 2148.160171297619
        "B2" = 5814.069708004545
        "E2" = 57548.65702072511
        "G2" = 8095.92571266189
        "I2" = 47220.76837855846
        "L2" = 63083.39540355118
        "M2" = 23203.20223601122
        "N2" = 13680.12751419989
        "O2" = 13142.91637740311
    },
    @{
        "A2" = 2148.160171297619
        "B2" = 5814.069708004545
        "E2" = 57548.65702072511
        "G2" = 8095.92571266189
        "I2" = 47220.76837855846
        "L2" = 63083.39540355118
        "M2" = 23203.20223601122
        "N2" = 13680.12751419989
        "O2" = 13142.91637740311
    }
)

for ($i = 0; $i -lt $sheetNames.Length; $i++) {
    $ws = $wb.Worksheets.Item($sheetNames[$i])
    $cellValues = $sheetValues[$i]
    foreach ($cellRef in $cellValues.Keys) {
        $ws.Range($cellRef).Value = $cellValues[$cellRef]
    }
}
